# Update "Project" interface mock data with additional properties: slug, description
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for the additional columns
$ws.Range("H1").Value = "slug"
$ws.Range("I1").Value = "description"

# Populate the new data for row 2, and fill in the previously-empty "images" cell
$ws.Range("G2").Value = "project-1-first.png,project-1-second.png,project-1-third.png,project-1-fourth.jpg,project-1-fifth.jpg,project-1-six.jpg"
$ws.Range("H2").Value = "project-one"
$ws.Range("I2").Value = "Cillum consectetur fugiat consequat cillum consequat exercitation qui dolore eu quis proident culpa duis. Exercitation exercitation non esse officia proident ipsum cupidatat nulla duis mollit laborum nostrud. Commodo enim exercitation veniam nisi ipsum ut laborum consectetur. Dolor  Dolor labore Lorem id cupidatat ea commodo incididunt incididunt elit ipsum labore Lorem veniam laboris. Amet adipisicing voluptate nostrud amet laborum sunt sint nulla est elit est dolor cupidatat nostrud. Velit aliqua fugiat enim aliqua cillum reprehenderit. Lorem adipisicing enim consectetur qui voluptate est eu nostrud cillum magna proident et in. Est elit culpa et anim consectetur dolore magna non incididunt veniam pariatur. Ut duis reprehenderit est aliquip eiusmod aliquip do magna esse ex. Minim proident enim Lorem dolor. Dolor labore Lorem id cupidatat ea commodo incididunt incididunt elit ipsum labore Lorem veniam laboris. Amet adipisicing voluptate nostrud amet laborum sunt sint nulla est elit est dolor cupidatat nostrud. Velit aliqua fugiat enim aliqua cillum reprehenderit. Lorem adipisicing enim consectetur qui voluptate est eu nostrud cillum magna proident et in. Est elit culpa et anim consectetur dolore magna non incididunt veniam pariatur. Ut duis reprehenderit est aliquip eiusmod aliquip do magna esse ex. Minim proident enim Lorem dolor.amet do velit nisi deserunt aliquip consequat. Minim aliquip dolore irure sint labore sunt incididunt qui nostrud."

# Adjust column widths for the newly populated / resized columns
$ws.Columns.Item(2).ColumnWidth = 22.5
$ws.Columns.Item(3).ColumnWidth = 47.833333333333336
$ws.Columns.Item(4).ColumnWidth = 27.5
$ws.Columns.Item(7).ColumnWidth = 35
$ws.Columns.Item(8).ColumnWidth = 15.166666666666668
$ws.Columns.Item(9).ColumnWidth = 18.666666666666664

# Row height for the data row
$ws.Rows.Item(2).RowHeight = 15

# Update view state: scroll the window over a bit and move the selection,
# matching where the editor's cursor ended up after adding the new columns
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("J11").Select()
